# Regenerate s_val data (filtered save games) for rows 2-9, columns B-E and G.
# Values are written as literal numbers (matching the source data, which has no formulas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = [double]"0.01253208636536152"
$ws.Range("C2").Value = [double]"6.194867796516235e-07"
$ws.Range("D2").Value = [double]"0.1496068669990043"
$ws.Range("E2").Value = [double]"0.5333859586016987"
$ws.Range("G2").Value = [double]"0.6955255314528441"

# Row 3
$ws.Range("B3").Value = [double]"1.174341637932841e-06"
$ws.Range("C3").Value = [double]"2.41451303395479e-12"
$ws.Range("D3").Value = [double]"0.7210945179870265"
$ws.Range("E3").Value = [double]"0.5333859586016987"
$ws.Range("G3").Value = [double]"1.254481650932778"

# Row 4
$ws.Range("B4").Value = [double]"1.445647641019636"
$ws.Range("C4").Value = [double]"1.626987699542094"
$ws.Range("D4").Value = [double]"3.223369029078222"
$ws.Range("E4").Value = [double]"13.86384647080068"
$ws.Range("G4").Value = [double]"20.15985084044064"

# Row 5
$ws.Range("B5").Value = [double]"3.272327238179451"
$ws.Range("C5").Value = [double]"1.626987699542094"
$ws.Range("D5").Value = [double]"3.223369029078222"
$ws.Range("E5").Value = [double]"0.5333859586016987"
$ws.Range("G5").Value = [double]"8.656069925401464"

# Row 6
$ws.Range("B6").Value = [double]"0.6545652718822623"
$ws.Range("C6").Value = [double]"0.04103571897497393"
$ws.Range("D6").Value = [double]"3.223369029078222"
$ws.Range("E6").Value = [double]"13.86384647080068"
$ws.Range("G6").Value = [double]"17.78281649073614"

# Row 7
$ws.Range("B7").Value = [double]"0.2881169905109251"
$ws.Range("C7").Value = [double]"0.04103571897497393"
$ws.Range("D7").Value = [double]"0.7210945179870265"
$ws.Range("E7").Value = [double]"0.5333859586016987"
$ws.Range("G7").Value = [double]"1.583633186074624"

# Row 8
$ws.Range("B8").Value = [double]"3.272327238179451"
$ws.Range("C8").Value = [double]"1.626987699542094"
$ws.Range("D8").Value = [double]"0.1496068669990043"
$ws.Range("E8").Value = [double]"0.5333859586016987"
$ws.Range("G8").Value = [double]"5.582307763322248"

# Row 9
$ws.Range("B9").Value = [double]"3.272327238179451"
$ws.Range("C9").Value = [double]"1.626987699542094"
$ws.Range("D9").Value = [double]"0.7210945179870265"
$ws.Range("E9").Value = [double]"13.86384647080068"
$ws.Range("G9").Value = [double]"19.48425592650926"
